# updated tilll 31st July entries
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Total Hours" row (currently row 81, with two blank spacer rows below it
# implicitly) needs to move down by 3 rows (to row 84) to make room for two
# new log entries (rows 79 and 80).
$null = $ws.Range("A81:G83").Insert(-4121)

# Clone the formatting of the last existing entry row (78) into the two new
# rows so number formats / alignment / wrap match the rest of the log.
$null = $ws.Range("A78:G78").Copy()
$null = $ws.Range("A79:G80").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 79 - new log entry (2022-07-27)
$ws.Range("A79").Value = 78
$ws.Range("B79").Value = 44769
$ws.Range("C79").Value = 0.21875
$ws.Range("D79").Value = 0.27083333333333331
$ws.Range("E79").Formula = "=D79-C79"
$ws.Range("F79").Value = "Code"
$ws.Range("G79").Value = "1. Self-Attention visualisation, Attention rollout attempt"

# Row 80 - new log entry (2022-07-27)
$ws.Range("A80").Value = 79
$ws.Range("B80").Value = 44769
$ws.Range("C80").Value = 0.21875
$ws.Range("D80").Value = 0.25
$ws.Range("E80").Formula = "=D80-C80"
$ws.Range("F80").Value = "Code"
$ws.Range("G80").Value = "1. Self-Attention visualisation, Attention rollout attempt"

# Move the view roughly to where the author left it (top visible row / active
# cell) - best effort, some of this is cosmetic only.
$excel.ActiveWindow.ScrollRow = 8
$null = $ws.Range("D88").Select()
